# Updated cryptos list (GitHub Actions refresh): refreshes Price (D) and
# Volume(1h) (E) columns for each coin row, and fixes the BNB/XRP row order
# (rows 5 and 6 had swapped names/links/prices).
#
# Numeric-looking price strings (e.g. "247.59") are forced to Text via
# NumberFormat "@" before assignment so Excel doesn't silently convert them
# to numbers (which would also truncate meaningful trailing zeros such as
# "73.70" -> 73.7). Prices that already contain thousand-separator dots
# (e.g. "35.781.28") are left alone since Excel already treats those as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.781.28'
$ws.Range('E2').Value = '  +0.30%  '

$ws.Range('D3').Value = '1.901.60'
$ws.Range('E3').Value = '  +0.06%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.59'
$ws.Range('E5').Value = '  -0.43%  '

$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.695'
$ws.Range('E6').Value = '  -0.02%  '

$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.28'
$ws.Range('E8').Value = '  -1.98%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '57.02'
$ws.Range('E9').Value = '  +9.72%  '

$ws.Range('E10').Value = '  +1.05%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0757'
$ws.Range('E11').Value = '  +1.72%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0985'
$ws.Range('E12').Value = '  +1.42%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.61'
$ws.Range('E13').Value = '  +11.11%  '

$ws.Range('E14').Value = '  +8.63%  '

$ws.Range('D15').Value = '2.171.57'
$ws.Range('E15').Value = '  -0.30%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.06'
$ws.Range('E16').Value = '  +1.90%  '

$ws.Range('D17').Value = '1.885.81'
$ws.Range('E17').Value = '  -0.74%  '

$ws.Range('D18').Value = '35.715.04'
$ws.Range('E18').Value = '  +0.13%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.70'
$ws.Range('E19').Value = '  -0.39%  '

$ws.Range('D20').Value = '0.0₃0833'
$ws.Range('E20').Value = '  +0.88%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '247.48'
$ws.Range('E21').Value = '  -0.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.09'
$ws.Range('E22').Value = '  +1.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.21'
$ws.Range('E23').Value = '  +4.55%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.70'
$ws.Range('E24').Value = '  +5.55%  '

$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -2.68%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.95'
$ws.Range('E27').Value = '  +0.70%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.73'
$ws.Range('E28').Value = '  +2.29%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.44'
$ws.Range('E29').Value = '  +0.02%  '

$ws.Range('E30').Value = '  +0.26%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.43'
$ws.Range('E31').Value = '  +3.74%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0610'
$ws.Range('E32').Value = '  +4.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.28'
$ws.Range('E33').Value = '  +0.75%  '

$ws.Range('E34').Value = '  -0.06%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.83'
$ws.Range('E35').Value = '  +18.12%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.49'
$ws.Range('E36').Value = '  -16.38%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.859'
$ws.Range('E37').Value = '  +0.25%  '

$ws.Range('E38').Value = '  +8.41%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.96'
$ws.Range('E39').Value = '  -2.93%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0230'
$ws.Range('E40').Value = '  +7.40%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.64'
$ws.Range('E41').Value = '  +1.41%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.07'
$ws.Range('E42').Value = '  -1.19%  '

$ws.Range('E43').Value = '  -0.46%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.36'
$ws.Range('E44').Value = '  +18.70%  '

$ws.Range('D45').Value = '1.318.17'
$ws.Range('E45').Value = '  +1.22%  '

$ws.Range('E46').Value = '  -0.87%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0812'
$ws.Range('E47').Value = '  +0.52%  '

$ws.Range('E48').Value = '  +0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.75'
$ws.Range('E49').Value = '  -0.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.42'
$ws.Range('E50').Value = '  +0.43%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.80'
$ws.Range('E51').Value = '  -1.84%  '
